# Refresh Universalis market-price snapshot columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) for the leves whose prices moved since the last
# scheduled pull. One crafting-job worksheet per tab (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Cells.Item(43, 8).Value = 1099
$ws.Cells.Item(43, 9).Value = 999
$ws.Cells.Item(43, 10).Value = 1149
$ws.Cells.Item(43, 11).Value = 999
$ws.Cells.Item(43, 12).Value = 1149
$ws.Cells.Item(43, 13).Value = -930
$ws.Cells.Item(43, 14).Value = -1287
# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Cells.Item(74, 8).Value = 5468.95
$ws.Cells.Item(74, 9).Value = 5404.625
$ws.Cells.Item(74, 11).Value = 5404.625
$ws.Cells.Item(74, 13).Value = -4468.625
# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Cells.Item(77, 8).Value = 5468.95
$ws.Cells.Item(77, 9).Value = 5404.625
$ws.Cells.Item(77, 11).Value = 27023.125
$ws.Cells.Item(77, 13).Value = -22343.125
# Row 114: Conserving Combat / Bluespirit Codex
$ws.Cells.Item(114, 8).Value = 99989.336
$ws.Cells.Item(114, 10).Value = 99989.336
$ws.Cells.Item(114, 12).Value = 99989.336
$ws.Cells.Item(114, 14).Value = -108667.336
# Row 118: Crafty Concoctions / Commanding Craftsman's Syrup
$ws.Cells.Item(118, 8).Value = 544.75
$ws.Cells.Item(118, 9).Value = 596.2857
$ws.Cells.Item(118, 11).Value = 1788.8571
$ws.Cells.Item(118, 13).Value = -131.8571000000002
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value = 1528.6052
$ws.Cells.Item(132, 9).Value = 1558.0555
$ws.Cells.Item(132, 11).Value = 4674.166499999999
$ws.Cells.Item(132, 13).Value = -2144.166499999999
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 331427.84
$ws.Cells.Item(137, 9).Value = 1740.6522
$ws.Cells.Item(137, 11).Value = 5221.9566
$ws.Cells.Item(137, 13).Value = -2671.9566

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Cells.Item(32, 8).Value = 5353.8696
$ws.Cells.Item(32, 9).Value = 2138.4038
$ws.Cells.Item(32, 11).Value = 2138.4038
$ws.Cells.Item(32, 13).Value = -1851.4038
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 2120.3333
$ws.Cells.Item(61, 10).Value = 2833.3333
$ws.Cells.Item(61, 12).Value = 2833.3333
$ws.Cells.Item(61, 14).Value = -3257.3333
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value = 2185.818
$ws.Cells.Item(74, 9).Value = 1511.6
$ws.Cells.Item(74, 10).Value = 3630.5715
$ws.Cells.Item(74, 11).Value = 1511.6
$ws.Cells.Item(74, 12).Value = 3630.5715
$ws.Cells.Item(74, 13).Value = -637.5999999999999
$ws.Cells.Item(74, 14).Value = -5378.5715
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value = 2185.818
$ws.Cells.Item(77, 9).Value = 1511.6
$ws.Cells.Item(77, 10).Value = 3630.5715
$ws.Cells.Item(77, 11).Value = 7558
$ws.Cells.Item(77, 12).Value = 18152.8575
$ws.Cells.Item(77, 13).Value = -3190
$ws.Cells.Item(77, 14).Value = -26888.8575
# Row 97: Ore for Me / High Steel Ingot
$ws.Cells.Item(97, 8).Value = 704
$ws.Cells.Item(97, 9).Value = 636.63635
$ws.Cells.Item(97, 11).Value = 636.63635
$ws.Cells.Item(97, 13).Value = -140.63635
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 2120.3333
$ws.Cells.Item(136, 10).Value = 2833.3333
$ws.Cells.Item(136, 12).Value = 8499.999899999999
$ws.Cells.Item(136, 14).Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector / Titanium Lump Hammer
$ws.Cells.Item(82, 8).Value = 9002.333000000001
$ws.Cells.Item(82, 9).Value = 9002.333000000001
$ws.Cells.Item(82, 11).Value = 9002.333000000001
$ws.Cells.Item(82, 13).Value = -8619.333000000001
# Row 85: The Clamor for Hammers (L) / Titanium Lump Hammer
$ws.Cells.Item(85, 8).Value = 9002.333000000001
$ws.Cells.Item(85, 9).Value = 9002.333000000001
$ws.Cells.Item(85, 11).Value = 9002.333000000001
$ws.Cells.Item(85, 13).Value = -7676.333000000001
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Cells.Item(86, 8).Value = 9567.444
$ws.Cells.Item(86, 9).Value = 6516.6665
$ws.Cells.Item(86, 11).Value = 6516.6665
$ws.Cells.Item(86, 13).Value = -5393.6665
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Cells.Item(89, 8).Value = 9567.444
$ws.Cells.Item(89, 9).Value = 6516.6665
$ws.Cells.Item(89, 11).Value = 32583.3325
$ws.Cells.Item(89, 13).Value = -26967.3325
# Row 97: File under Dull / High Steel File
$ws.Cells.Item(97, 8).Value = 6825
$ws.Cells.Item(97, 9).Value = 6825
$ws.Cells.Item(97, 11).Value = 6825
$ws.Cells.Item(97, 13).Value = -5834
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 2231.0715
$ws.Cells.Item(107, 9).Value = 1946.7
$ws.Cells.Item(107, 11).Value = 1946.7
$ws.Cells.Item(107, 13).Value = -26.70000000000005

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 3361.1538
$ws.Cells.Item(31, 9).Value = 2185
$ws.Cells.Item(31, 11).Value = 2185
$ws.Cells.Item(31, 13).Value = -1890
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 3361.1538
$ws.Cells.Item(34, 9).Value = 2185
$ws.Cells.Item(34, 11).Value = 2185
$ws.Cells.Item(34, 13).Value = -1983
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 1290.1945
$ws.Cells.Item(58, 9).Value = 1131.5927
$ws.Cells.Item(58, 11).Value = 1131.5927
$ws.Cells.Item(58, 13).Value = -928.5926999999999
# Row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws.Cells.Item(80, 8).Value = 15000
$ws.Cells.Item(80, 10).Value = 15000
$ws.Cells.Item(80, 12).Value = 15000
$ws.Cells.Item(80, 14).Value = -17246
# Row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws.Cells.Item(83, 8).Value = 15000
$ws.Cells.Item(83, 10).Value = 15000
$ws.Cells.Item(83, 12).Value = 45000
$ws.Cells.Item(83, 14).Value = -56232
# Row 99: O Pine Pine / Lumber
$ws.Cells.Item(99, 8).Value = 5293349
$ws.Cells.Item(99, 9).Value = 7938534
$ws.Cells.Item(99, 10).Value = 2978.4285
$ws.Cells.Item(99, 11).Value = 7938534
$ws.Cells.Item(99, 12).Value = 2978.4285
$ws.Cells.Item(99, 13).Value = -7937036
$ws.Cells.Item(99, 14).Value = -5974.4285
# Row 126: A Better Conductor / Red Pine Lumber
$ws.Cells.Item(126, 8).Value = 5293349
$ws.Cells.Item(126, 9).Value = 7938534
$ws.Cells.Item(126, 10).Value = 2978.4285
$ws.Cells.Item(126, 11).Value = 23815602
$ws.Cells.Item(126, 12).Value = 8935.2855
$ws.Cells.Item(126, 13).Value = -23813132
$ws.Cells.Item(126, 14).Value = -13875.2855
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 1494.6154
$ws.Cells.Item(132, 9).Value = 1348.3334
$ws.Cells.Item(132, 11).Value = 4045.0002
$ws.Cells.Item(132, 13).Value = -1515.0002
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Cells.Item(134, 8).Value = 4938.25
$ws.Cells.Item(134, 9).Value = 5501
$ws.Cells.Item(134, 11).Value = 16503
$ws.Cells.Item(134, 13).Value = -13968
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 1290.1945
$ws.Cells.Item(136, 9).Value = 1131.5927
$ws.Cells.Item(136, 11).Value = 3394.7781
$ws.Cells.Item(136, 13).Value = -844.7780999999995

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Cells.Item(2, 8).Value = 29.705883
$ws.Cells.Item(2, 10).Value = 50.833332
$ws.Cells.Item(2, 12).Value = 304.999992
$ws.Cells.Item(2, 14).Value = -530.999992
# Row 7: It's Always Sunny in Vylbrand / Raisins
$ws.Cells.Item(7, 8).Value = 7290.421
$ws.Cells.Item(7, 9).Value = 5769.2144
$ws.Cells.Item(7, 10).Value = 11549.8
$ws.Cells.Item(7, 11).Value = 17307.6432
$ws.Cells.Item(7, 12).Value = 34649.39999999999
$ws.Cells.Item(7, 13).Value = -17195.6432
$ws.Cells.Item(7, 14).Value = -34873.39999999999
# Row 18: Fisher of Men / Salt Cod
$ws.Cells.Item(18, 8).Value = 9349.333000000001
$ws.Cells.Item(18, 9).Value = 10819.2
$ws.Cells.Item(18, 11).Value = 32457.6
$ws.Cells.Item(18, 13).Value = -32288.6
# Row 132: More Mezcal / Cooking Mezcal
$ws.Cells.Item(132, 8).Value = 5123.905
$ws.Cells.Item(132, 9).Value = 1162.8572
$ws.Cells.Item(132, 11).Value = 10465.7148
$ws.Cells.Item(132, 13).Value = -7935.7148
# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Cells.Item(140, 8).Value = 1459
$ws.Cells.Item(140, 9).Value = 1136.7368
$ws.Cells.Item(140, 11).Value = 3410.2104
$ws.Cells.Item(140, 13).Value = 1769.7896

$ws = $wb.Worksheets.Item("GSM")
# Row 57: Gold Is So Last Year / Electrum Circlet (Amber)
$ws.Cells.Item(57, 8).Value = 23916.5
$ws.Cells.Item(57, 9).Value = 5499.5
$ws.Cells.Item(57, 11).Value = 5499.5
$ws.Cells.Item(57, 13).Value = -4679.5
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Cells.Item(122, 8).Value = 1120267.6
$ws.Cells.Item(122, 9).Value = 1437343
$ws.Cells.Item(122, 11).Value = 4312029
$ws.Cells.Item(122, 13).Value = -4309579
# Row 132: On Board for Lar / Lar Ingot
$ws.Cells.Item(132, 8).Value = 3853.5334
$ws.Cells.Item(132, 9).Value = 2896.7222
$ws.Cells.Item(132, 11).Value = 8690.1666
$ws.Cells.Item(132, 13).Value = -6160.1666

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Cells.Item(22, 8).Value = 1661.3334
$ws.Cells.Item(22, 9).Value = 1487.5
$ws.Cells.Item(22, 10).Value = 1800.4
$ws.Cells.Item(22, 11).Value = 1487.5
$ws.Cells.Item(22, 12).Value = 1800.4
$ws.Cells.Item(22, 13).Value = -1192.5
$ws.Cells.Item(22, 14).Value = -2390.4
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Cells.Item(27, 8).Value = 1661.3334
$ws.Cells.Item(27, 9).Value = 1487.5
$ws.Cells.Item(27, 10).Value = 1800.4
$ws.Cells.Item(27, 11).Value = 1487.5
$ws.Cells.Item(27, 12).Value = 1800.4
$ws.Cells.Item(27, 13).Value = -1380.5
$ws.Cells.Item(27, 14).Value = -2014.4
# Row 40: Best Served Toad / Toad Leather
$ws.Cells.Item(40, 8).Value = 7409919.5
$ws.Cells.Item(40, 9).Value = 2776.3157
$ws.Cells.Item(40, 11).Value = 2776.3157
$ws.Cells.Item(40, 13).Value = -2640.3157
# Row 46: Supply Side Logic / Boar Leather
$ws.Cells.Item(46, 8).Value = 3594.25
$ws.Cells.Item(46, 9).Value = 1974
$ws.Cells.Item(46, 10).Value = 3774.2778
$ws.Cells.Item(46, 11).Value = 1974
$ws.Cells.Item(46, 12).Value = 3774.2778
$ws.Cells.Item(46, 13).Value = -1786
$ws.Cells.Item(46, 14).Value = -4150.2778
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Cells.Item(55, 8).Value = 7118.75
$ws.Cells.Item(55, 9).Value = 5987.5
$ws.Cells.Item(55, 10).Value = 8250
$ws.Cells.Item(55, 11).Value = 5987.5
$ws.Cells.Item(55, 12).Value = 8250
$ws.Cells.Item(55, 13).Value = -5814.5
$ws.Cells.Item(55, 14).Value = -8596
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Cells.Item(82, 8).Value = 2399.7778
$ws.Cells.Item(82, 9).Value = 2471.1428
$ws.Cells.Item(82, 11).Value = 2471.1428
$ws.Cells.Item(82, 13).Value = -2110.1428
# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Cells.Item(85, 8).Value = 2399.7778
$ws.Cells.Item(85, 9).Value = 2471.1428
$ws.Cells.Item(85, 11).Value = 2471.1428
$ws.Cells.Item(85, 13).Value = -1223.1428
# Row 122: Hell on Leather / Gaja Leather
$ws.Cells.Item(122, 8).Value = 46157330
$ws.Cells.Item(122, 10).Value = 16670497
$ws.Cells.Item(122, 12).Value = 50011491
$ws.Cells.Item(122, 14).Value = -50016391
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Cells.Item(132, 8).Value = 2898.75
$ws.Cells.Item(132, 9).Value = 2678.7
$ws.Cells.Item(132, 11).Value = 8036.099999999999
$ws.Cells.Item(132, 13).Value = -5506.099999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 950.5454999999999
$ws.Cells.Item(132, 9).Value = 802.19354
$ws.Cells.Item(132, 11).Value = 2406.58062
$ws.Cells.Item(132, 13).Value = 123.4193800000003
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 3388.8
$ws.Cells.Item(136, 9).Value = 2898
$ws.Cells.Item(136, 11).Value = 8694
$ws.Cells.Item(136, 13).Value = -6144
